$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F24, F25: convert from text status to numeric percentage completion
$ws.Range("F24").Value = 0.8
$ws.Range("F24").NumberFormat = "0%"

$ws.Range("F25").Value = 0.3
$ws.Range("F25").NumberFormat = "0%"

# E27/F27, E28/F28: fill in member + status ("Hoàn Thành")
$ws.Range("E27").Value = "Cả nhóm"
$ws.Range("F27").Value = "Hoàn Thành"

$ws.Range("E28").Value = "Cả nhóm"
$ws.Range("F28").Value = "Hoàn Thành"

# E29/F29: member + numeric percentage completion
$ws.Range("E29").Value = "Cả nhóm"
$ws.Range("F29").Value = 0.7
$ws.Range("F29").NumberFormat = "0%"

# Update sheet view scroll/selection to match new work area
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
